# "implementing formik in request section"
# - insert a new blank row before row 25 (pushes old rows 25-45 down to 26-46,
#   extends the filler-row block by one, dimension grows to A1:V46)
# - put the new "* строка для FCA" note into (still-blank) A24
# - split the "Приемка" defined name into "Приемка1" (A23, unchanged) and a
#   new "Приемка2" (A24, the FCA line)
# - re-point "ВСД" / "ВСД_далее" at their shifted rows (26 / 27)
# - grow the Print_Area to match the new last row (46)
# - update the active selection to the cell the author ended up on (D19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift rows 25+ down by one, creating a fresh blank row 25 -------------
$ws.Rows("25:25").Insert()

# --- the row that stayed at 24 gets the new note ---------------------------
$ws.Range("A24").Value = "* строка для FCA"

# --- defined names -----------------------------------------------------
# "Приемка" (-> $A$23) becomes "Приемка1" (same target cell)
$wb.Names.Item("Приемка").Name = "Приемка1"

# add the new "Приемка2" name pointing at the new FCA row ($A$24).
# Names.Add() chokes on non-ASCII name arguments in this host, so mint it
# with an ASCII placeholder and rename it right after.
$wb.Names.Add("TempName_Priemka2", "=Request_Contract!`$A`$24")
$wb.Names.Item("TempName_Priemka2").Name = "Приемка2"

# "ВСД" / "ВСД_далее" followed their rows down by one
$wb.Names.Item("ВСД").RefersTo = "=Request_Contract!`$A`$26"
$wb.Names.Item("ВСД_далее").RefersTo = "=Request_Contract!`$A`$27"

# --- print area grows by the inserted row -----------------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$G`$46"

# --- restore the author's final selection ------------------------------
$ws.Range("D19").Select()
